$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring the formatting (date/time number formats etc.) for the new rows in
# line with the existing rows by copying row 43's formats down, the same
# way a user would continue the log by dragging/filling the row.
$ws.Range("A43:E43").Copy()
$ws.Range("A44:E44").PasteSpecial(-4122)
$ws.Range("A43:E43").Copy()
$ws.Range("A45:E45").PasteSpecial(-4122)
$ws.Range("A43:E43").Copy()
$ws.Range("A46:E46").PasteSpecial(-4122)

# Row 44 - 10/28/2021
$ws.Cells.Item(44, 1).Value = 44497
$ws.Cells.Item(44, 2).Value = 0.3833333333333333
$ws.Cells.Item(44, 3).Value = 0.4145833333333333
$ws.Cells.Item(44, 4).Formula = "=C44-B44"
$ws.Cells.Item(44, 5).Value = "worked on jshs form"

# Row 45 - 10/29/2021
$ws.Cells.Item(45, 1).Value = 44498
$ws.Cells.Item(45, 2).Value = 0.3833333333333333
$ws.Cells.Item(45, 3).Value = 0.4145833333333333
$ws.Cells.Item(45, 4).Formula = "=C45-B45"
$ws.Cells.Item(45, 5).Value = "worked on presentation"

# Row 46 - 11/1/2021
$ws.Cells.Item(46, 1).Value = 44501
$ws.Cells.Item(46, 2).Value = 0.3833333333333333
$ws.Cells.Item(46, 3).Value = 0.4145833333333333
$ws.Cells.Item(46, 4).Formula = "=C46-B46"
$ws.Cells.Item(46, 5).Value = "presented"

$ws.Range("A47").Select()
